$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the NERC vocabulary URI to the d18O attribute definition (row 16, column B)
$ws.Range("B16").Value = "Enrichment of oxygen-18 in dissolved oxygen {18O in O2 CAS 14797-71-8} {delta(18)O} in the water body by mass spectrometry URI http://vocab.nerc.ac.uk/collection/P01/current/D18OMXDG/"

# Leave the selection on the cell that was last edited
$ws.Range("B16").Select()
